# Update the "want to go" (想去人数) counts for the 南宁-漫展信息 workbook.
# These values live in column F on both the "展览" and "全部类型" sheets
# (which mirror the same data).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 658
    $ws.Range("F3").Value = 3906
    $ws.Range("F4").Value = 109
}
